$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 39 (pushes existing rows 39+ down to 41+,
# and auto-shifts formulas / defined names that reference rows below
# the insertion point).
$ws.Rows("39:40").Insert()

# New row 39 -> "keT1": derived from keT3_ (the drug's T-cell elimination rate)
$ws.Range("E39").Value = "keT1"
$ws.Range("F39").Formula = "=keT3_"
$ws.Range("G39").Value = "1/d"
$ws.Range("H39").Value = "calc"

# New row 40 -> "keDT1": derived from row above (F39)
$ws.Range("E40").Value = "keDT1"
$ws.Range("F40").Formula = "=F39"
$ws.Range("G40").Value = "1/d"
$ws.Range("H40").Value = "calc"

# The defined name "keT3_" previously pointed at F39 (now occupied by the
# new "keT1" row); repoint it to F41, where the original keT3 formula
# (=LN(2)/1) now lives after the insert shifted it down.
$wb.Names.Item("keT3_").RefersTo = "=Sheet1!`$F`$41"

# Update the view so the active selection matches the edited area.
$ws.Range("G41").Select()
